# Update the dSF column (F) with recalculated values following a data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -1
    5  = -4
    6  = -2
    7  = 2
    8  = -1
    9  = 1
    10 = 4
    11 = -1
    12 = -1
    13 = -2
    14 = -1
    15 = 1
    16 = 7
    18 = 6
    19 = -2
    20 = -2
    22 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
